$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V1").Value = "HO_chg"

$values = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = -35.946843853820596
    8  = 0
    9  = 0
    10 = -37.291246100201867
    11 = 0
    12 = 0
    13 = 0
    14 = -37.738771295818275
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = -5.6249999999999982
    28 = 0
    29 = 0
    30 = -17.668414683340053
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($row, 22)
    $cell.Value = $values[$row]
    if ($values[$row] -eq 0) {
        $cell.NumberFormat = "0.00"
    }
}

$ws.Range("X34").Select()
